$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "ASYNC"
$ws.Range("C2").Value = "Done"

$ws.Columns.Item(2).ColumnWidth = 77.5

$ws.Range("C4").Select() | Out-Null
